$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Sheet1"

# Header row
$ws.Range("A1").Value = "Nome"
$ws.Range("B1").Value = "Faltas"
$ws.Range("C1").Value = "Matemática"
$ws.Range("D1").Value = "Ciências"
$ws.Range("E1").Value = "Física"

# Row 2 - Vinicius
$ws.Range("A2").Value = "Vinicius"
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 7
$ws.Range("D2").Value = "5,5"
$ws.Range("E2").Value = 8

# Row 3 - Marco
$ws.Range("A3").Value = "Marco"
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = "8,5"
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = 5

# Row 4 - Vivian
$ws.Range("A4").Value = "Vivian"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 9
$ws.Range("E4").Value = 7

# Row 5 - Ana
$ws.Range("A5").Value = "Ana"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "7,5"
$ws.Range("D5").Value = "7,5"
$ws.Range("E5").Value = 9

# Row 6 - Tatiana
$ws.Range("A6").Value = "Tatiana"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 4

# ---- Formatting ----
# Main style (font Arial 10 black, thin gray border, vertical bottom):
# header row (A1:E1), columns A,B and E for data rows (2-6)
$grayColor = 10132122

$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Name = "Arial"
$headerRange.Font.Size = 10
$headerRange.Font.Color = 0
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Color = $grayColor
$headerRange.VerticalAlignment = -4107

$colA = $ws.Range("A2:A6")
$colA.Font.Name = "Arial"
$colA.Font.Size = 10
$colA.Font.Color = 0
$colA.Borders.LineStyle = 1
$colA.Borders.Color = $grayColor
$colA.VerticalAlignment = -4107

$colB = $ws.Range("B2:B6")
$colB.Font.Name = "Arial"
$colB.Font.Size = 10
$colB.Font.Color = 0
$colB.Borders.LineStyle = 1
$colB.Borders.Color = $grayColor
$colB.VerticalAlignment = -4107

$colE = $ws.Range("E2:E6")
$colE.Font.Name = "Arial"
$colE.Font.Size = 10
$colE.Font.Color = 0
$colE.Borders.LineStyle = 1
$colE.Borders.Color = $grayColor
$colE.VerticalAlignment = -4107

# Right-aligned style (font Arial 10 black, thin gray border, vertical bottom, horizontal right):
# columns C,D for data rows (2-6)
$colC = $ws.Range("C2:C6")
$colC.Font.Name = "Arial"
$colC.Font.Size = 10
$colC.Font.Color = 0
$colC.Borders.LineStyle = 1
$colC.Borders.Color = $grayColor
$colC.VerticalAlignment = -4107
$colC.HorizontalAlignment = -4152

$colD = $ws.Range("D2:D6")
$colD.Font.Name = "Arial"
$colD.Font.Size = 10
$colD.Font.Color = 0
$colD.Borders.LineStyle = 1
$colD.Borders.Color = $grayColor
$colD.VerticalAlignment = -4107
$colD.HorizontalAlignment = -4152

# Empty formatted rows 7-10 (A:E): thin gray border, theme Arial 10 font, vertical bottom, no value
$emptyRange = $ws.Range("A7:E10")
$emptyRange.Font.Name = "Arial"
$emptyRange.Font.Size = 10
$emptyRange.Borders.LineStyle = 1
$emptyRange.Borders.Color = $grayColor
$emptyRange.VerticalAlignment = -4107
